{"js": "// Apply the documented edits to the Science Gateway Integration Coordination task doc.\n// 1) Extend the Summary sentence about when coordination can begin.\n// 2) Drop \"RP\" before \"contacts:\" in the RP-contacts intro sentence.\n// 3) Swap \"the RP\" -> \"the Gateway Administrator\" and drop \"RP\" before \"staff\" later\n//    in the same paragraph.\n// 4) Drop \"RP\" before \"contacts\" in the \"enter and update\" sentence.\n// 5) Rename the hyperlink display text for the contacts spreadsheet.\n// 6) Drop \"RP\" -> \"gateway\" in the \"staff changes\" sentence and remove the trailing\n//    blank paragraph plus the \"ACCESS Integration Roadmaps...\" paragraph that followed it.\n\nconst body = context.document.body;\n\n// 1) Summary paragraph - append the new sentence.\nconst summaryHits = body.search(\n  \"A science gateway wishing to integrate with ACCESS contacts ACCESS to start the integration process, provides basic science gateway information, and identifies science gateway staff contacts that will be involved in coordination, technical integration, and ongoing support activities. \",\n  { matchCase: true }\n);\nsummaryHits.load(\"items\");\nawait context.sync();\nsummaryHits.items[0].insertText(\n  \"A science gateway wishing to integrate with ACCESS contacts ACCESS to start the integration process, provides basic science gateway information, and identifies science gateway staff contacts that will be involved in coordination, technical integration, and ongoing support activities. This coordination can begin before or after an ACCESS allocation is awarded. \",\n  Word.InsertLocation.replace\n);\n\n// 2) \"... must provide each of the following RP contacts:\" -> drop \" RP\"\nconst mustProvideHits = body.search(\n  \"he Science Gateway Integration Coordinator must provide each of the following RP contacts:\",\n  { matchCase: true }\n);\nmustProvideHits.load(\"items\");\nawait context.sync();\nmustProvideHits.items[0].insertText(\n  \"he Science Gateway Integration Coordinator must provide each of the following contacts:\",\n  Word.InsertLocation.replace\n);\n\n// 3) \"If the RP doesn't know...\" -> \"If the Gateway Administrator doesn't know...\"\n//    and \"...add other RP staff replacing...\" -> \"...add other staff replacing...\"\nconst placeholderHits = body.search(\n  \"If the RP doesn\\u2019t know who some of these contacts will be, or wishes not to identify them until relevant integration effort ramps up, please name the Integration Coordinator as the placeholder for those roles. In other words, ACCESS needs each someone to be the contact in each of the above areas, even if it\\u2019s the coordinator who will eventually hand off effort to someone else. As integration activities ramps up the Integration Coordinator can add other RP staff replacing themselves as necessary.\",\n  { matchCase: true }\n);\nplaceholderHits.load(\"items\");\nawait context.sync();\nplaceholderHits.items[0].insertText(\n  \"If the Gateway Administrator doesn\\u2019t know who some of these contacts will be, or wishes not to identify them until relevant integration effort ramps up, please name the Integration Coordinator as the placeholder for those roles. In other words, ACCESS needs each someone to be the contact in each of the above areas, even if it\\u2019s the coordinator who will eventually hand off effort to someone else. As integration activities ramps up the Integration Coordinator can add other staff replacing themselves as necessary.\",\n  Word.InsertLocation.replace\n);\n\n// 4) \"... enter and update their RP contacts in the resource specific tab here:\" -> drop \" RP\"\nconst enterUpdateHits = body.search(\n  \"The Integration Coordinator should enter and update their RP contacts in the resource specific tab here:\",\n  { matchCase: true }\n);\nenterUpdateHits.load(\"items\");\nawait context.sync();\nenterUpdateHits.items[0].insertText(\n  \"The Integration Coordinator should enter and update their contacts in the resource specific tab here:\",\n  Word.InsertLocation.replace\n);\n\n// 5) Hyperlink display text rename.\nconst hyperlinkHits = body.search(\"ACCESS Allocated Resource Contacts\", { matchCase: true });\nhyperlinkHits.load(\"items\");\nawait context.sync();\nhyperlinkHits.items[0].insertText(\"ACCESS Science Gateways Contacts\", Word.InsertLocation.replace);\n\n// 6) \"...updating when RP staff changes....\" -> \"...updating when gateway staff changes....\"\nconst maintainingHits = body.search(\n  \"The Integration Coordinator is responsible for maintaining accurate contact information in this spreadsheet and updating when RP staff changes. We recommend that the Integration Coordinator review and correct their contacts every 6 months. This sub-task should take ~1 hour to complete initially, and ~\\u00bd hour annually to keep up-to-date.\",\n  { matchCase: true }\n);\nmaintainingHits.load(\"items\");\nawait context.sync();\nmaintainingHits.items[0].insertText(\n  \"The Integration Coordinator is responsible for maintaining accurate contact information in this spreadsheet and updating when gateway staff changes. We recommend that the Integration Coordinator review and correct their contacts every 6 months. This sub-task should take ~1 hour to complete initially, and ~\\u00bd hour annually to keep up-to-date.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Remove the now-empty paragraph and the trailing \"ACCESS Integration Roadmaps...\" paragraph\n// that used to sit right after the sentence edited above.\nconst roadmapHits = body.search(\n  \"ACCESS Integration Roadmaps task lists which RP staff contacts normally perform that task.\",\n  { matchCase: true }\n);\nroadmapHits.load(\"items\");\nawait context.sync();\n\nconst roadmapParagraph = roadmapHits.items[0].paragraphs.getFirst();\nconst blankParagraph = roadmapParagraph.getPrevious();\nroadmapParagraph.delete();\nblankParagraph.delete();\nawait context.sync();\n", "ps1": "# Apply the documented edits to the Science Gateway Integration Coordination task doc.\n# 1) Extend the Summary sentence about when coordination can begin.\n# 2) Drop \"RP\" before \"contacts:\" in the RP-contacts intro sentence.\n# 3) Swap \"the RP\" -> \"the Gateway Administrator\" and drop \"RP\" before \"staff\" later\n#    in the same paragraph.\n# 4) Drop \"RP\" before \"contacts\" in the \"enter and update\" sentence.\n# 5) Rename the hyperlink display text for the contacts spreadsheet.\n# 6) Drop \"RP\" -> \"gateway\" in the \"staff changes\" sentence and remove the trailing\n#    blank paragraph plus the \"ACCESS Integration Roadmaps...\" paragraph that followed it.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($searchText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $searchText\n  $find.Replacement.Text = $replaceText\n  $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Summary paragraph - append the new sentence.\nReplace-DocText `\n  \"A science gateway wishing to integrate with ACCESS contacts ACCESS to start the integration process, provides basic science gateway information, and identifies science gateway staff contacts that will be involved in coordination, technical integration, and ongoing support activities. \" `\n  \"A science gateway wishing to integrate with ACCESS contacts ACCESS to start the integration process, provides basic science gateway information, and identifies science gateway staff contacts that will be involved in coordination, technical integration, and ongoing support activities. This coordination can begin before or after an ACCESS allocation is awarded. \"\n\n# 2) \"... must provide each of the following RP contacts:\" -> drop \" RP\"\nReplace-DocText `\n  \"he Science Gateway Integration Coordinator must provide each of the following RP contacts:\" `\n  \"he Science Gateway Integration Coordinator must provide each of the following contacts:\"\n\n# 3) \"If the RP doesn't know...\" -> \"If the Gateway Administrator doesn't know...\"\n#    and \"...add other RP staff replacing...\" -> \"...add other staff replacing...\"\nReplace-DocText `\n  \"If the RP doesn\u2019t know who some of these contacts will be, or wishes not to identify them until relevant integration effort ramps up, please name the Integration Coordinator as the placeholder for those roles. In other words, ACCESS needs each someone to be the contact in each of the above areas, even if it\u2019s the coordinator who will eventually hand off effort to someone else. As integration activities ramps up the Integration Coordinator can add other RP staff replacing themselves as necessary.\" `\n  \"If the Gateway Administrator doesn\u2019t know who some of these contacts will be, or wishes not to identify them until relevant integration effort ramps up, please name the Integration Coordinator as the placeholder for those roles. In other words, ACCESS needs each someone to be the contact in each of the above areas, even if it\u2019s the coordinator who will eventually hand off effort to someone else. As integration activities ramps up the Integration Coordinator can add other staff replacing themselves as necessary.\"\n\n# 4) \"... enter and update their RP contacts in the resource specific tab here:\" -> drop \" RP\"\nReplace-DocText `\n  \"The Integration Coordinator should enter and update their RP contacts in the resource specific tab here:\" `\n  \"The Integration Coordinator should enter and update their contacts in the resource specific tab here:\"\n\n# 5) Hyperlink display text rename.\nReplace-DocText `\n  \"ACCESS Allocated Resource Contacts\" `\n  \"ACCESS Science Gateways Contacts\"\n\n# 6) \"...updating when RP staff changes....\" -> \"...updating when gateway staff changes....\"\nReplace-DocText `\n  \"The Integration Coordinator is responsible for maintaining accurate contact information in this spreadsheet and updating when RP staff changes. We recommend that the Integration Coordinator review and correct their contacts every 6 months. This sub-task should take ~1 hour to complete initially, and ~\u00bd hour annually to keep up-to-date.\" `\n  \"The Integration Coordinator is responsible for maintaining accurate contact information in this spreadsheet and updating when gateway staff changes. We recommend that the Integration Coordinator review and correct their contacts every 6 months. This sub-task should take ~1 hour to complete initially, and ~\u00bd hour annually to keep up-to-date.\"\n\n# Remove the now stale \"ACCESS Integration Roadmaps...\" paragraph and the blank\n# paragraph before it, walking paragraphs from the end so indices stay valid.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"ACCESS Integration Roadmaps task lists which RP staff contacts normally perform that task.\") {\n    $blank = $d.Paragraphs.Item($i - 1)\n    $p.Range.Delete()\n    $blank.Range.Delete()\n    break\n  }\n}\n"}
